$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2 through 15 from serial date 45204 to 45205
for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
